# Refresh cryptos list (Price / Volume(1h) columns) with freshly scraped values.
# Price-column (D) values are free-form scraped text (e.g. "37.424.05", trailing
# zeros like "19.00") that must stay literal text, not get reinterpreted as a
# number -- so those cells are forced to Text format for the write, then the
# style is reset back to Normal (matching the original unstyled cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.424.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.051.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.28"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.387"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0792"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("E11").Value = "  -1.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.353.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.760"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.28%  "
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.053.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.297.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.27%  "
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0831"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("E25").Value = "  -3.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("E28").Value = "  -5.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.59%  "
$ws.Range("E30").Value = "  -3.08%  "
$ws.Range("E31").Value = "  -2.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0617"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("E36").Value = "  +1.56%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("E38").Value = "  -4.88%  "
$ws.Range("E39").Value = "  -1.74%  "
$ws.Range("E40").Value = "  -4.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.477.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0942"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.28%  "
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("E47").Value = "  -4.14%  "
$ws.Range("E48").Value = "  -4.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("E50").Value = "  -2.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.241.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.48%  "
